$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last updated" banner date (A1)
$ws.Range("A1").Value = "Last updated: November 19, 2020"

# Fix the "databse" -> "database" typo in the cdc_report_dt description (F4)
$ws.Range("F4").Value = "This date was populated using the date at which a case record was first submitted to the database. If missing, then the report date entered on the case report form was used. If missing, then the date at which the case first appeared in the database was used."

# Drop the "(MM/DD/YYYY)" suffix from the pos_spec_dt description (A12)
$ws.Range("A12").Value = "Date of first positive specimen collection"

# Move the active selection to F5
$ws.Range("F5").Select()
